$wb = $excel.ActiveWorkbook
$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# The handback transform failed for the e692daac-dab8-4262-b19e-44c5c2cee0b7 file
# (row 7 in every sheet). Its shared "Status" text moves from
# "Ready for handoff" to "Handback transform failed" everywhere it is shown:
# Overview!E7 (zh-cn status), Overview!F7 (de-de status), zh-cn!C7 and de-de!C7.
$wsOverview.Cells.Item(7, 5).Value = "Handback transform failed"
$wsOverview.Cells.Item(7, 6).Value = "Handback transform failed"
$wsZhCn.Cells.Item(7, 3).Value = "Handback transform failed"
$wsDeDe.Cells.Item(7, 3).Value = "Handback transform failed"

# Populate the Error Detail column (P, the 16th column) with the handback
# mismatch message for that same row, and widen the column to fit it.
$wsZhCn.Cells.Item(7, 16).Value = "Handback file name: jnbxxrdl.ow1 is different with handoff file name: e692daac-dab8-4262-b19e-44c5c2cee0b7.7a2513fccf940db10006b55693c6d3176fa7a21c.zh-cn."
$wsDeDe.Cells.Item(7, 16).Value = "Handback file name: jnbxxrdl.ow1 is different with handoff file name: e692daac-dab8-4262-b19e-44c5c2cee0b7.7a2513fccf940db10006b55693c6d3176fa7a21c.de-de."

# ColumnWidth is in character units; Excel stores the column's <col width=.../>
# 5/6 of a character wider than the ColumnWidth value, so back that off to
# land on a stored width of exactly 40.
$wsZhCn.Columns.Item(16).ColumnWidth = 39.166666666666664
$wsDeDe.Columns.Item(16).ColumnWidth = 39.166666666666664
